$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append, matching the source diff (dates are Excel serial
# numbers continuing the existing series, formatted the same way as column A
# in the preceding rows).
$data = @(
    @(239, 44313, 1, 10, 107.7121930202499),
    @(240, 44314, 0, 10, 107.7121930202499),
    @(241, 44315, 6, 14, 150.7970702283498),
    @(242, 44316, 8, 21, 226.1956053425248),
    @(243, 44317, 6, 26, 280.0517018526497),
    @(244, 44318, 5, 28, 301.5941404566997)
)

# Template row for copying the date cell's style/number format (col A).
$lastRow = 238
$srcA = $ws.Range("A" + $lastRow)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $b = $row[2]
    $c = $row[3]
    $d = $row[4]

    $srcA.Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = $dateSerial
    $ws.Range("B" + $r).Value = $b
    $ws.Range("C" + $r).Value = $c
    $ws.Range("D" + $r).Value = $d
}

$excel.CutCopyMode = 0
